# Updates the crypto price/volume table to the latest snapshot values.
# (GitHub Actions scheduled refresh of cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (e.g. "69.281.08", "1.00", "0.160")
# that must stay text. Assigning a numeric-looking string via .Value would
# otherwise be auto-converted to a real number (losing formatting like
# trailing zeros), so each target cell is temporarily switched to the Text
# number format, written, and then has that temporary formatting cleared so
# the cell keeps its original (unstyled) look while the stored type is text.
$priceTextCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D10", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D23", "D25", "D26", "D28", "D30", "D33", "D34", "D36", "D37", "D38", "D40", "D43", "D45", "D46", "D47", "D48", "D50")
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Write the refreshed price + volume(1h) figures.
# (Row 45/46 also swap: Monero now ranks above TheGraph.)
$ws.Range("D2").Value = "69.281.08"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "3.776.26"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "622.59"
$ws.Range("E5").Value = "  +3.98%  "
$ws.Range("D6").Value = "165.48"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("D7").Value = "3.774.84"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "35.68"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").Value = "4.407.90"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "3.728.72"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "69.253.51"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "17.66"
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("D19").Value = "7.10"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").Value = "468.18"
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "0.703"
$ws.Range("E24").Value = "  +4.97%  "
$ws.Range("D25").Value = "83.32"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "12.01"
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("E27").Value = "  +3.65%  "
$ws.Range("D28").Value = "10.03"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "3.924.79"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +4.19%  "
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "7.31"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").Value = "28.82"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "3.726.46"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "8.99"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").Value = "0.160"
$ws.Range("E38").Value = "  +13.13%  "
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("D40").Value = "3.39"
$ws.Range("E40").Value = "  +6.97%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "154.57"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.299"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").Value = "43.21"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "46.81"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").Value = "8.42"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("E51").Value = "  +0.32%  "

# Drop the temporary Text formatting so the cells end up unstyled, matching
# the rest of the sheet.
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).ClearFormats()
}
